$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-17 03:20:32"
$wsZh.Range("H2").Value = "2016-03-17 03:21:12"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-17 03:20:40"
$wsDe.Range("H2").Value = "2016-03-17 03:21:25"
